# Nudge three shapes inside the nested "组合 1" group left, on both slides.
#   id=14 (直接连接符 13, straight connector) : new x =  2948354 EMU (was 3024554)
#   id=19 (文本框 18, Approve/Disapprove)      : new x =  2435200 EMU (was 2501875)
#   id=30 (文本框 29, Match Orders)            : new x =  2199943 EMU (was 2266618)
#
# Shape.Left/.Top in this object model are slide-space coordinates expressed
# in points and stored as 32-bit floats, while the OOXML offsets are EMU
# (1 pt = 12700 EMU). Both ancestor groups here have ext == chExt (1:1
# scale), so an EMU delta in the shape's own <a:off> equals the same EMU
# delta in slide space. We still target the *absolute* new slide-space
# value (rather than Left - delta) and add half an EMU (in points) so the
# float32 round-trip/floor lands exactly on the intended EMU instead of
# one unit short.

$p = $ppt.ActivePresentation

$emuPerPoint = 12700
$halfEmuInPoints = 0.5 / $emuPerPoint

# New absolute x offsets (EMU, slide space) keyed by shape Id.
$newXEmu = @{
    14 = 2948354
    19 = 2435200
    30 = 2199943
}

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    $topGroup = $s.Shapes.Item(1)

    foreach ($shape in $topGroup.GroupItems) {
        if ($newXEmu.ContainsKey($shape.Id)) {
            $targetEmu = $newXEmu[$shape.Id]
            $shape.Left = ($targetEmu / $emuPerPoint) + $halfEmuInPoints
        }
    }
}
